# Insert a new data row above row 253 (this pushes the existing rows
# 253..365 down to 254..366, matching the diff's shift pattern), then
# populate the new row 253 with the new "Berenjena" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(253).Insert()

$ws.Range("A253").Value = 10
$ws.Range("B253").Value = "Vega Modelo de Temuco"
$ws.Range("C253").Value = "La Araucanía"
$ws.Range("D253").Value = 44875
$ws.Range("E253").Value = 9
$ws.Range("F253").Value = 100112001
$ws.Range("G253").Value = "Berenjena"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 55
$ws.Range("K253").Value = 15000
$ws.Range("L253").Value = 15000
$ws.Range("M253").Value = 15000
$ws.Range("N253").Value = "$/caja 40 unidades"
$ws.Range("O253").Value = "Región de Arica y Parinacota"
$ws.Range("P253").Value = 375
$ws.Range("Q253").Value = 40
$ws.Range("R253").Value = "Hortaliza"
